$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("POS")

# Remove the standalone title row ("Avnet Japan March 2020") at the top of
# the sheet so the header row becomes row 1 and all data shifts up by one.
$ws.Rows.Item(1).Delete()
